$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.397.74'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '2.646.67'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = "'597.47"
$ws.Range('E5').Value = '  -0.13%  '
$ws.Range('D6').Value = "'158.77"
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -1.21%  '
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('E10').Value = '  -1.13%  '
$ws.Range('D11').Value = "'5.29"
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').Value = "'0.350"
$ws.Range('E12').Value = '  -1.28%  '
$ws.Range('D13').Value = "'28.01"
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').Value = '3.124.27'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('E15').Value = '  -2.89%  '
$ws.Range('D16').Value = '68.247.59'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '2.637.80'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = "'11.41"
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').Value = "'363.54"
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('D20').Value = "'7.45"
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').Value = "'4.40"
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = "'4.77"
$ws.Range('E22').Value = '  -2.70%  '
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = "'74.37"
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = "'9.77"
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('E28').Value = '  -2.93%  '
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('D30').Value = "'560.21"
$ws.Range('E30').Value = '  -2.56%  '
$ws.Range('D31').Value = "'8.05"
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('E32').Value = '  -2.01%  '
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').Value = "'1.66"
$ws.Range('E34').Value = '  +4.46%  '
$ws.Range('E35').Value = '  -1.63%  '
$ws.Range('D36').Value = "'1.00"
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').Value = "'160.43"
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('D38').Value = "'19.64"
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('E39').Value = '  -1.23%  '
$ws.Range('D41').Value = "'5.32"
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('D42').Value = "'2.61"
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').Value = '0.0₆0321'
$ws.Range('E43').Value = '  -4.69%  '
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').Value = "'158.08"
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = "'3.81"
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('D47').Value = "'22.02"
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('E48').Value = '  -1.54%  '
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('D50').Value = "'0.574"
$ws.Range('E50').Value = '  +1.08%  '
$ws.Range('E51').Value = '  -1.04%  '
